$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp + value rows 2-27 with the new dataset
$data = @(
    @(2, 45786, 4.479, 1.174),
    @(3, 45786.01041666666, 12.146, 0),
    @(4, 45786.02083333334, 3.927, 0.216),
    @(5, 45786.03125, 0.438, 16.108),
    @(6, 45786.04166666666, 0, 70.155),
    @(7, 45786.05208333334, 0, 62.613),
    @(8, 45786.0625, 0, 35.32),
    @(9, 45786.07291666666, 0, 19.555),
    @(10, 45786.08333333334, 0.259, 16.731),
    @(11, 45786.09375, 27.505, 0),
    @(12, 45786.10416666666, 16.479, 0),
    @(13, 45786.11458333334, 6.513, 0.735),
    @(14, 45786.125, 0, 8.257999999999999),
    @(15, 45786.13541666666, 0, 18.763),
    @(16, 45786.14583333334, 0, 19.018),
    @(17, 45786.15625, 0, 20.37),
    @(18, 45786.16666666666, 0, 16.178),
    @(19, 45786.17708333334, 0, 22.429),
    @(20, 45786.1875, 0, 28.596),
    @(21, 45786.19791666666, 0, 44.248),
    @(22, 45786.20833333334, 0.538, 11.246),
    @(23, 45786.21875, 0, 21.588),
    @(24, 45786.22916666666, 0.243, 10.626),
    @(25, 45786.23958333334, 8.673999999999999, 0),
    @(26, 45786.25, 2.683, 15.72),
    @(27, 45786.26041666666, 0.611, 1.419)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Remove the now-unused trailing rows (28-40) so the sheet dimension shrinks to A1:C27
$ws.Rows("28:40").Delete() | Out-Null

